# Apply updated transition-probability matrix values to Sheet1 (B2:S19)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(0.1949458483754513,0.5451263537906137,0,0,0,0,0,0,0.02888086642599278,0,0,0,0,0,0.1444043321299639,0,0,0.08664259927797834),
    @(0.01176470588235294,0.05294117647058823,0,0,0,0,0,0,0.01764705882352941,0,0,0,0,0,0.7823529411764706,0,0,0.1352941176470588),
    @(0,0,0,0,0,0,0,0,0.09615384615384616,0,0,0,0,0,0.6346153846153846,0,0,0.2692307692307692),
    @(0,0,0,0,0,0,0,0,0.3333333333333333,0,0,0,0,0,0.6666666666666666,0,0,0),
    @(0.07983193277310924,0,0.01680672268907563,0.004201680672268907,0.07142857142857142,0,0,0,0.2605042016806723,0,0,0,0,0.01260504201680672,0,0.1680672268907563,0.04621848739495799,0.3403361344537815),
    @(0.1005291005291005,0,0.01587301587301587,0,0.0582010582010582,0,0,0,0.1322751322751323,0,0,0,0,0.01587301587301587,0,0.201058201058201,0.08994708994708994,0.3862433862433862),
    @(0.08893709327548807,0,0.02169197396963124,0.002169197396963124,0.06941431670281996,0,0,0,0.123644251626898,0,0,0,0,0.02603036876355748,0,0.2017353579175705,0.1301518438177874,0.3362255965292842),
    @(0.05617977528089887,0,0.01123595505617977,0,0.0449438202247191,0,0,0,0.1235955056179775,0,0,0,0,0.02247191011235955,0,0.1460674157303371,0.1460674157303371,0.449438202247191),
    @(0.1018444266238974,0,0.0264635124298316,0.0008019246190858059,0.07056936647955092,0,0,0,0.1234963913392141,0,0,0,0,0.01443464314354451,0,0.2510024057738572,0.09783480352846832,0.3135525260625501),
    @(0,0,0,0,0.006802721088435374,0.1394557823129252,0,0,0.1054421768707483,0.1972789115646258,0.5272108843537415,0,0,0,0,0,0,0.02380952380952381),
    @(0,0,0,0,0,0.8129032258064516,0,0,0.1483870967741935,0.006451612903225806,0.01935483870967742,0,0,0,0,0,0,0.01290322580645161),
    @(0,0,0,0,0,0.6923076923076923,0,0,0.2564102564102564,0,0,0,0,0,0,0,0,0.05128205128205128),
    @(0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0.03539823008849557,0,0.1504424778761062,0.02654867256637168,0.4026548672566372,0.06637168141592921,0,0.02654867256637168,0.004424778761061947,0.04424778761061947,0,0,0,0.2433628318584071),
    @(0,0,0,0,0.03589743589743589,0,0.1487179487179487,0.07179487179487179,0.3897435897435897,0.1230769230769231,0,0.01025641025641026,0,0.07179487179487179,0,0,0,0.1487179487179487),
    @(0,0,0,0,0.02772277227722772,0,0.1742574257425742,0.04752475247524753,0.4396039603960396,0.1168316831683168,0,0.01386138613861386,0,0.06732673267326733,0,0,0,0.1128712871287129),
    @(0,0,0,0,0.02232142857142857,0,0.2142857142857143,0.01339285714285714,0.4151785714285715,0.09375,0,0.01339285714285714,0.004464285714285714,0.08482142857142858,0,0,0,0.1383928571428572),
    @(0,0,0,0,0.02312673450508788,0,0.2395929694727104,0.03700277520814061,0.3718778908418132,0.1073080481036078,0,0.01942645698427382,0,0.07770582793709528,0,0,0,0.123959296947271)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value2 = $rowVals[$c]
    }
}
